$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V2").Value = 74.76470624280951
$ws.Range("W2").Value = 87.38143681982832
$ws.Range("Y2").Value = 398.135471261432
$ws.Range("Z2").Value = 465.3218246424835
$ws.Range("AB2").Value = 87.70875919968069
$ws.Range("AC2").Value = 75.04476757467819
$ws.Range("V3").Value = 132.2260097199337
$ws.Range("W3").Value = 141.4163303275339
$ws.Range("Y3").Value = 704.1272190904584
$ws.Range("Z3").Value = 753.0673247904316
$ws.Range("AC3").Value = 132.7213155166232
$ws.Range("V4").Value = 122.1976242665746
$ws.Range("Y4").Value = 650.7242677634285
$ws.Range("AC4").Value = 122.6553647048522
$ws.Range("V5").Value = 131.3789820895892
$ws.Range("Y5").Value = 699.6166450278322
$ws.Range("AC5").Value = 131.8711150029997
$ws.Range("V6").Value = 97.83623074704535
$ws.Range("Y6").Value = 520.995477577541
$ws.Range("AC6").Value = 98.20271577005931
$ws.Range("V7").Value = 138.3698287686405
$ws.Range("Y7").Value = 736.8441575394352
$ws.Range("AC7").Value = 138.8881487150806
$ws.Range("V8").Value = 77.26073897646847
$ws.Range("Y8").Value = 411.4272932806598
$ws.Range("AC8").Value = 77.5501502046573
$ws.Range("V9").Value = 124.2640596500611
$ws.Range("Y9").Value = 661.7284068362659
$ws.Range("AC9").Value = 124.7295407546879
$ws.Range("V10").Value = 128.653324736496
$ws.Range("Y10").Value = 685.1020307224317
$ws.Range("AC10").Value = 129.1352476020526
$ws.Range("V11").Value = 128.0286754852751
$ws.Range("Y11").Value = 681.7756614165692
$ws.Range("AC11").Value = 128.508258475374
$ws.Range("V12").Value = 151.1919147648076
$ws.Range("Y12").Value = 805.124065361979
$ws.Range("AC12").Value = 151.7582650006968
$ws.Range("V13").Value = 125.1488976079331
$ws.Range("Y13").Value = 666.4403276750007
$ws.Range("AC13").Value = 125.6176932296553
$ws.Range("V14").Value = 102.8037492387761
$ws.Range("Y14").Value = 547.4484045680127
$ws.Range("AC14").Value = 103.1888421038427
$ws.Range("V15").Value = 83.73573264114532
$ws.Range("Y15").Value = 445.907795962349
$ws.Range("AC15").Value = 84.04939856704782
$ws.Range("V16").Value = 117.9593776088418
$ws.Range("Y16").Value = 628.1548440982227
$ws.Range("AC16").Value = 118.4012419865636
$ws.Range("V17").Value = 109.584754508139
$ws.Range("Y17").Value = 583.5584739338465
$ws.Range("AC17").Value = 109.9952483606846
$ws.Range("V18").Value = 149.0554807383192
$ws.Range("W18").Value = 157.201814133566
$ws.Range("Y18").Value = 793.7471709595254
$ws.Range("Z18").Value = 837.1278574941061
$ws.Range("AB18").Value = 157.790676869085
$ws.Range("AC18").Value = 149.6138280997377
$ws.Range("V19").Value = 117.2700155373638
$ws.Range("Y19").Value = 624.4838674169756
$ws.Range("AC19").Value = 117.709297631685
$ws.Range("V20").Value = 101.016643147634
$ws.Range("Y20").Value = 537.9317440801041
$ws.Range("AC20").Value = 101.3950416867649
$ws.Range("V21").Value = 91.83728812922089
$ws.Range("Y21").Value = 489.0500321094459
$ws.Range("AC21").Value = 92.18130169553048
$ws.Range("V22").Value = 142.6150474863387
$ws.Range("W22").Value = 156.2302863745529
$ws.Range("Y22").Value = 759.4507086745327
$ws.Range("Z22").Value = 831.9542979147767
$ws.Range("AB22").Value = 156.8155098613969
$ws.Range("AC22").Value = 143.1492696099945
$ws.Range("V23").Value = 114.1672449182469
$ws.Range("W23").Value = 124.8707620997791
$ws.Range("Y23").Value = 607.9610573273278
$ws.Range("Z23").Value = 664.9592062050144
$ws.Range("AB23").Value = 125.3385158528875
$ws.Range("AC23").Value = 114.5949043350276
$ws.Range("V24").Value = 109.5847545081391
$ws.Range("Y24").Value = 583.5584739338469
$ws.Range("AC24").Value = 109.9952483606847
$ws.Range("V25").Value = 107.007531963467
$ws.Range("Y25").Value = 569.8343016079873
$ws.Range("AC25").Value = 107.4083717905418
$ws.Range("V26").Value = 138.3117967366042
$ws.Range("Y26").Value = 736.5351265596571
$ws.Range("AC26").Value = 138.8298993006857
$ws.Range("V27").Value = 117.1313970136564
$ws.Range("Y27").Value = 623.7456989142789
$ws.Range("AC27").Value = 117.5701598564441
$ws.Range("V28").Value = 138.5388472684103
$ws.Range("Y28").Value = 737.744211367497
$ws.Range("AC28").Value = 139.0578003417432
$ws.Range("V29").Value = 144.0202353673912
$ws.Range("Y29").Value = 766.9335861891824
$ws.Range("AC29").Value = 144.5597211884421
$ws.Range("V30").Value = 122.1440286444172
$ws.Range("Y30").Value = 650.4388614620131
$ws.Range("AC30").Value = 122.601568318697
$ws.Range("V31").Value = 139.8265451330602
$ws.Range("Y31").Value = 744.6014334706605
$ws.Range("AC31").Value = 140.350321797592
$ws.Range("V32").Value = 116.601605581278
$ws.Range("Y32").Value = 620.9244645083653
$ws.Range("AC32").Value = 117.038383876789
$ws.Range("V33").Value = 129.7936318496911
$ws.Range("Y33").Value = 691.1743706366724
$ws.Range("AC33").Value = 130.2798261949994
$ws.Range("V34").Value = 128.1993362163265
$ws.Range("Y34").Value = 682.6844604207715
$ws.Range("AC34").Value = 128.6795584849571
$ws.Range("V35").Value = 140.6183461625234
$ws.Range("Y35").Value = 748.8179159776178
$ws.Range("AC35").Value = 141.1450888368481
$ws.Range("V36").Value = 110.8262623198196
$ws.Range("Y36").Value = 590.1697257198557
$ws.Range("AC36").Value = 111.2414067401095
$ws.Range("V37").Value = 149.4562851418366
$ws.Range("Y37").Value = 795.8815263003974
$ws.Range("AC37").Value = 150.0161338776429
$ws.Range("V38").Value = 143.9453970158842
$ws.Range("Y38").Value = 766.5350585436796
$ws.Range("AC38").Value = 144.4846024997352
$ws.Range("V39").Value = 100.5152574056096
$ws.Range("Y39").Value = 535.2617750704439
$ws.Range("AC39").Value = 100.8917778024225
$ws.Range("V40").Value = 126.679122586973
$ws.Range("Y40").Value = 674.5890501642925
$ws.Range("AC40").Value = 127.1536502829215
$ws.Range("V41").Value = 164.928113212401
$ws.Range("Y41").Value = 878.2717859523898
$ws.Range("AC41").Value = 165.5459179142458
$ws.Range("V42").Value = 140.5643479812306
$ws.Range("Y42").Value = 748.5303659765998
$ws.Range("AC42").Value = 141.0908883836101
$ws.Range("V43").Value = 148.5047810309958
$ws.Range("Y43").Value = 790.8145962392197
$ws.Range("AC43").Value = 149.0610655247693
$ws.Range("V44").Value = 109.3538629453479
$ws.Range("Y44").Value = 582.3289349470473
$ws.Range("AC44").Value = 109.7634919005128
$ws.Range("V45").Value = 134.9675463389499
$ws.Range("Y45").Value = 718.7263933351568
$ws.Range("AC45").Value = 135.4731216656819
$ws.Range("V46").Value = 156.0132407271587
$ws.Range("Y46").Value = 830.798490911003
$ws.Range("AC46").Value = 156.5976511820771
$ws.Range("V47").Value = 125.3260165936287
$ws.Range("Y47").Value = 667.3835180435984
$ws.Range("AC47").Value = 125.7954756858776
$ws.Range("V48").Value = 133.8208017725088
$ws.Range("Y48").Value = 712.6197728276979
$ws.Range("AC48").Value = 134.3220815054142
$ws.Range("V49").Value = 115.7232480753053
$ws.Range("Y49").Value = 616.2470532383908
$ws.Range("AC49").Value = 116.1567361288648
$ws.Range("V50").Value = 91.75567240425681
$ws.Range("Y50").Value = 488.6154137346288
$ws.Range("AC50").Value = 92.09938024598344
$ws.Range("V51").Value = 133.9333992494458
$ws.Range("Y51").Value = 713.2193745889564
$ws.Range("AC51").Value = 134.4351007615695
$ws.Range("V52").Value = 115.3309391323641
$ws.Range("Y52").Value = 614.1579377489149
$ws.Range("AC52").Value = 115.7629576347063
$ws.Range("V53").Value = 110.6468593319974
$ws.Range("Y53").Value = 589.2143726302529
$ws.Range("AC53").Value = 111.0613317261103
$ws.Range("V54").Value = 120.0275454120225
$ws.Range("Y54").Value = 639.1682086166754
$ws.Range("AC54").Value = 120.4771569455695
$ws.Range("V55").Value = 75.2293353344853
$ws.Range("Y55").Value = 400.6097045149587
$ws.Range("AC55").Value = 75.51113712182794
$ws.Range("V56").Value = 147.1640940312401
$ws.Range("Y56").Value = 783.6751974869755
$ws.Range("AC56").Value = 147.7153564416585
$ws.Range("V57").Value = 134.8997502706081
$ws.Range("W57").Value = 151.4513975732123
$ws.Range("Y57").Value = 718.3653671106804
$ws.Range("Z57").Value = 806.5058578600721
$ws.Range("AB57").Value = 152.0187198065128
$ws.Range("AC57").Value = 135.4050716398494
$ws.Range("V58").Value = 117.1573194811519
$ws.Range("Y58").Value = 623.8837407033963
$ws.Range("AC58").Value = 117.5961794269865
$ws.Range("V59").Value = 120.491413973953
$ws.Range("Y59").Value = 641.6383919128922
$ws.Range("AC59").Value = 120.9427631141033
$ws.Range("V60").Value = 129.9771983471983
$ws.Range("Y60").Value = 692.1518951621549
$ws.Range("AC60").Value = 130.4640803147868
$ws.Range("V61").Value = 138.0784969123962
$ws.Range("Y61").Value = 735.2927631488444
$ws.Range("AC61").Value = 138.5957255579837
$ws.Range("V62").Value = 123.0907255124977
$ws.Range("Y62").Value = 655.4801920932184
$ws.Range("AC62").Value = 123.5518114213456
$ws.Range("V63").Value = 74.09656429725177
$ws.Range("Y63").Value = 394.5774955569565
$ws.Range("AC63").Value = 74.37412283425166
$ws.Range("V64").Value = 143.1672248336066
$ws.Range("Y64").Value = 762.391152092587
$ws.Range("AC64").Value = 143.7035153600033
$ws.Range("V65").Value = 145.6748989141434
$ws.Range("W65").Value = 155.2819605154788
$ws.Range("Y65").Value = 775.7449663720431
$ws.Range("Z65").Value = 826.9042926143387
$ws.Range("AC65").Value = 146.2205829442155
$ws.Range("V66").Value = 126.2588893153503
$ws.Range("W66").Value = 128.0406825025471
$ws.Range("Y66").Value = 672.3512326158107
$ws.Range("Z66").Value = 681.8396009372368
$ws.Range("AB66").Value = 128.5203104697671
$ws.Range("AC66").Value = 126.731842858257
$ws.Range("V67").Value = 120.6066597640275
$ws.Range("Y67").Value = 642.2520964166349
$ws.Range("AC67").Value = 121.0584406037201
$ws.Range("V68").Value = 126.1175473065487
$ws.Range("Y68").Value = 671.5985610664769
$ws.Range("AC68").Value = 126.5899713960111
$ws.Range("V69").Value = 127.6272003643167
$ws.Range("W69").Value = 142.2391787457397
$ws.Range("Y69").Value = 679.6377343850179
$ws.Range("Z69").Value = 757.4491402114023
$ws.Range("AC69").Value = 128.1052794675853
$ws.Range("V70").Value = 147.3349073440194
$ws.Range("Y70").Value = 784.584809016248
$ws.Range("AC70").Value = 147.886809604526
$ws.Range("V71").Value = 153.3497314802519
$ws.Range("Y71").Value = 816.6148264185282
$ws.Range("AC71").Value = 153.9241646880891
$ws.Range("V72").Value = 118.5080041052441
$ws.Range("Y72").Value = 631.0763785985023
$ws.Range("AC72").Value = 118.951923584564
$ws.Range("V73").Value = 124.7332851746996
$ws.Range("W73").Value = 136.1972742453382
$ws.Range("Y73").Value = 664.2271169197797
$ws.Range("Z73").Value = 725.2749150125272
$ws.Range("AB73").Value = 136.7074560134329
$ws.Range("AC73").Value = 125.2005239525921
$ws.Range("V74").Value = 115.3859311676189
$ws.Range("Y74").Value = 614.4507802005483
$ws.Range("AC74").Value = 115.8181556647869
$ws.Range("V75").Value = 122.2897305109332
$ws.Range("W75").Value = 133.9148791591713
$ws.Range("Y75").Value = 651.2147500357024
$ws.Range("Z75").Value = 713.1207517863028
$ws.Range("AB75").Value = 134.4165112968347
$ws.Range("AC75").Value = 122.7478159702611
$ws.Range("V76").Value = 87.87514035515903
$ws.Range("Y76").Value = 434.5258213432918
$ws.Range("AC76").Value = 88.20431210052891
$ws.Range("V77").Value = 82.48950412206295
$ws.Range("Y77").Value = 407.8948765939116
$ws.Range("AC77").Value = 82.79850179690948
$ws.Range("V78").Value = 101.1635607966164
$ws.Range("Y78").Value = 500.2345278482509
$ws.Range("AC78").Value = 101.5425096749984
$ws.Range("V79").Value = 73.41168272055538
$ws.Range("Y79").Value = 363.0067798630799
$ws.Range("AC79").Value = 73.68667575765323
$ws.Range("V80").Value = 138.2100790573799
$ws.Range("Y80").Value = 735.9934616727244
$ws.Range("AC80").Value = 138.7278005969092
$ws.Range("V81").Value = 87.80798481579551
$ws.Range("Y81").Value = 434.19375
$ws.Range("AC81").Value = 88.13690500303399
$ws.Range("V82").Value = 116.1816522044278
$ws.Range("Y82").Value = 574.4961276319125
$ws.Range("AC82").Value = 116.6168573953558
$ws.Range("V83").Value = 71.716867548991
$ws.Range("Y83").Value = 354.6262418466096
$ws.Range("AC83").Value = 71.98551197297884
$ws.Range("V84").Value = 98.71902205402468
$ws.Range("Y84").Value = 525.6964996331652
$ws.Range("AC84").Value = 99.08881392757841
$ws.Range("V85").Value = 115.4861820167965
$ws.Range("W85").Value = 139.0822775386679
$ws.Range("Y85").Value = 614.9846339543777
$ws.Range("Z85").Value = 740.6380750315111
$ws.Range("AC85").Value = 115.9187820439137
$ws.Range("V86").Value = 120.0046548321027
$ws.Range("Y86").Value = 639.0463121727271
$ws.Range("AC86").Value = 120.4541806197594
$ws.Range("V87").Value = 120.7854112625108
$ws.Range("Y87").Value = 643.2039802086504
$ws.Range("AC87").Value = 121.2378616879646
$ws.Range("V88").Value = 126.1072760740915
$ws.Range("W88").Value = 143.8144671274709
$ws.Range("Y88").Value = 671.5438649112953
$ws.Range("Z88").Value = 765.8378334030316
$ws.Range("AB88").Value = 144.3531821606698
$ws.Range("AC88").Value = 126.5796616885142
$ws.Range("V89").Value = 107.4991176294036
$ws.Range("Y89").Value = 572.4520834546438
$ws.Range("AC89").Value = 107.9017988886628
$ws.Range("V90").Value = 146.7083957683407
$ws.Range("Y90").Value = 781.2485224986028
$ws.Range("AC90").Value = 147.2579511773028
$ws.Range("V91").Value = 121.9552858362728
$ws.Range("Y91").Value = 649.4337721539144
$ws.Range("AC91").Value = 122.4121184983153
$ws.Range("V92").Value = 137.4396186972558
$ws.Range("W92").Value = 158.2134146233427
$ws.Range("Y92").Value = 731.890622057866
$ws.Range("Z92").Value = 842.5148115525827
$ws.Range("AB92").Value = 158.8060667160962
$ws.Range("AC92").Value = 137.9544541670679
$ws.Range("V93").Value = 128.0691390838696
$ws.Range("Y93").Value = 681.9911373370257
$ws.Range("AC93").Value = 128.5488736466808
$ws.Range("V94").Value = 97.87870520780282
$ws.Range("Y94").Value = 521.2216616997033
$ws.Range("AC94").Value = 98.24534933602364
$ws.Range("V95").Value = 143.0549963582182
$ws.Range("Y95").Value = 761.793515330762
$ws.Range("AC95").Value = 143.5908664876404
$ws.Range("V96").Value = 156.8344146973914
$ws.Range("Y96").Value = 835.1713895961713
$ws.Range("AC96").Value = 157.4219011902875
$ws.Range("V97").Value = 109.8650103276247
$ws.Range("Y97").Value = 585.0508864419936
$ws.Range("AC97").Value = 110.2765539912645
$ws.Range("V98").Value = 134.1524042530606
$ws.Range("Y98").Value = 714.3856155160533
$ws.Range("AC98").Value = 134.654926136668
$ws.Range("V99").Value = 110.8401756520711
$ws.Range("W99").Value = 136.1371500065007
$ws.Range("Y99").Value = 590.2438167097243
$ws.Range("Z99").Value = 724.9547426562543
$ws.Range("AB99").Value = 136.6471065550332
$ws.Range("AC99").Value = 111.255372190353
$ws.Range("V100").Value = 116.9212132125896
$ws.Range("Y100").Value = 622.6264324730058
$ws.Range("AC100").Value = 117.359188727262
$ws.Range("V101").Value = 143.106761695206
$ws.Range("Y101").Value = 762.069175035352
$ws.Range("AC101").Value = 143.6428257325555
$ws.Range("V102").Value = 103.6569133119592
$ws.Range("Y102").Value = 472.0981270896576
$ws.Range("Z102").Value = 521.9985140176923
$ws.Range("AC102").Value = 104.0452020468238
$ws.Range("V103").Value = 103.6569133119592
$ws.Range("Y103").Value = 472.0981270896576
$ws.Range("Z103").Value = 521.9985140176923
$ws.Range("AC103").Value = 104.0452020468238
$ws.Range("V104").Value = 103.6569133119592
$ws.Range("Y104").Value = 472.0981270896576
$ws.Range("Z104").Value = 521.9985140176923
$ws.Range("AC104").Value = 104.0452020468238
$ws.Range("V105").Value = 87.85419708945554
$ws.Range("Y105").Value = 400.1257666054041
$ws.Range("AC105").Value = 88.18329038338516
$ws.Range("V106").Value = 90.36413564431933
$ws.Range("Y106").Value = 411.557105365176
$ws.Range("AC106").Value = 90.70263092442542
$ws.Range("V107").Value = 155.8168351747937
$ws.Range("W107").Value = 166.6254805025663
$ws.Range("Y107").Value = 709.6568256250691
$ws.Range("Z107").Value = 758.8840411824239
$ws.Range("AB107").Value = 167.2496433775096
$ws.Range("AC107").Value = 156.4005099135787
$ws.Range("V108").Value = 66.58783364358995
$ws.Range("Y108").Value = 303.2696088054325
$ws.Range("AC108").Value = 66.83726520446452
$ws.Range("V109").Value = 113.9648335894923
$ws.Range("Y109").Value = 519.0448255946325
$ws.Range("AC109").Value = 114.3917347930844
$ws.Range("V110").Value = 96.67990004295457
$ws.Range("Y110").Value = 806.3724721902718
$ws.Range("AC110").Value = 97.04205356340077
$ws.Range("V111").Value = 129.6150774221181
$ws.Range("Y111").Value = 1081.073008635402
$ws.Range("AC111").Value = 130.1006029198738
$ws.Range("V113").Value = 95.44570053327124
$ws.Range("Y113").Value = 796.0784554467998
$ws.Range("AC113").Value = 95.80323086216298
$ws.Range("V115").Value = 93.64908218816682
$ws.Range("Y115").Value = 781.0935043258265
$ws.Range("AC115").Value = 93.9998825591429
$ws.Range("V116").Value = 117.9593223258183
$ws.Range("Y116").Value = 983.8565236362276
$ws.Range("AC116").Value = 118.4011864964553
$ws.Range("V121").Value = 111.3493450454741
$ws.Range("Y121").Value = 963.5522569710803
$ws.Range("AC121").Value = 111.7664488828766
$ws.Range("V122").Value = 93.38022670743987
$ws.Range("Y122").Value = 808.0579922915497
$ws.Range("AC122").Value = 93.73001997188408
$ws.Range("V127").Value = 108.1215154419398
$ws.Range("Y127").Value = 1011.481624972125
$ws.Range("AC127").Value = 108.5265281429861
$ws.Range("V128").Value = 107.1706226525133
$ws.Range("Y128").Value = 963.5241888504934
$ws.Range("AC128").Value = 107.5720734014774
$ws.Range("V129").Value = 148.6171602139234
$ws.Range("Y129").Value = 1336.151878194078
$ws.Range("AC129").Value = 149.1738656692071
$ws.Range("V130").Value = 124.079768523649
$ws.Range("Y130").Value = 943.924319331269
$ws.Range("AC130").Value = 124.5445592916061
$ws.Range("V131").Value = 92.23047934723739
$ws.Range("Y131").Value = 701.6340655313586
$ws.Range("AC131").Value = 92.57596576968098
$ws.Range("V132").Value = 144.3172704632135
$ws.Range("Y132").Value = 1097.879073362164
$ws.Range("AC132").Value = 144.857868948899
$ws.Range("V133").Value = 120.7694491830507
$ws.Range("Y133").Value = 918.7413989605902
$ws.Range("AC133").Value = 121.2218398161043
